# Assignment 2 is completed
# Adds a "Message Length" column (L) to the walletHubPostReviews sheet,
# converts J2 to a text value "4", and adds L2 = "200" (also stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell L1: "Message Length" ---
# Start from the same look as the other header cells (border + centered,
# style index 2) by copying A1's formatting, then set the text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = "Message Length"

# --- J2: change the numeric 4 into the text value "4" ---
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "4"

# --- L2: new cell holding the text value "200" ---
# Copy A1's border/center formatting first so the new cell has the same
# bordered look as the rest of row 2, then switch it to Text format.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "200"

# --- Column width for the new column L ---
$ws.Range("L1").ColumnWidth = 13.29

# --- Update selection to the newly active cell ---
$ws.Range("L2").Select() | Out-Null

Write-Host "Assignment 2 edits applied"
